# Auto-applies updated "want-to-go count" (F) and "lowest price" (G) values
# to the 展览 (sheet1), 演出 (sheet2), and 全部类型 (sheet4) worksheets,
# matching a refreshed scrape (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) updates: column F ---
$wsExpo.Range("F3").Value = 521
$wsExpo.Range("F4").Value = 267
$wsExpo.Range("F5").Value = 479
$wsExpo.Range("F6").Value = 1120
$wsExpo.Range("F8").Value = 23
$wsExpo.Range("F9").Value = 107
$wsExpo.Range("F10").Value = 107
$wsExpo.Range("F11").Value = 1115
$wsExpo.Range("F13").Value = 91
$wsExpo.Range("F14").Value = 759
$wsExpo.Range("F15").Value = 796
$wsExpo.Range("F16").Value = 178
$wsExpo.Range("F17").Value = 34
$wsExpo.Range("F18").Value = 60
$wsExpo.Range("F19").Value = 663
$wsExpo.Range("F20").Value = 162
$wsExpo.Range("F21").Value = 1706
$wsExpo.Range("F22").Value = 2143
$wsExpo.Range("F23").Value = 578
$wsExpo.Range("F25").Value = 1825
$wsExpo.Range("F26").Value = 282
$wsExpo.Range("F27").Value = 2660
$wsExpo.Range("F28").Value = 478
$wsExpo.Range("F30").Value = 667
$wsExpo.Range("F34").Value = 921
$wsExpo.Range("F35").Value = 1639
$wsExpo.Range("F36").Value = 304
$wsExpo.Range("F38").Value = 523
$wsExpo.Range("F39").Value = 132
$wsExpo.Range("F40").Value = 108
$wsExpo.Range("F41").Value = 146

# --- 演出 (sheet2) updates: column G ---
$wsShow.Range("G3").Value = 144

# --- 全部类型 (sheet4) updates: columns F and G ---
$wsAll.Range("F4").Value = 521
$wsAll.Range("F5").Value = 267
$wsAll.Range("F6").Value = 479
$wsAll.Range("F7").Value = 1120
$wsAll.Range("F9").Value = 23
$wsAll.Range("F10").Value = 107
$wsAll.Range("F11").Value = 107
$wsAll.Range("F12").Value = 1115
$wsAll.Range("F14").Value = 759
$wsAll.Range("F16").Value = 178
$wsAll.Range("G17").Value = 144
$wsAll.Range("G18").Value = 144
$wsAll.Range("F20").Value = 34
$wsAll.Range("F22").Value = 60
$wsAll.Range("F23").Value = 663
$wsAll.Range("F24").Value = 162
$wsAll.Range("F25").Value = 1706
$wsAll.Range("F26").Value = 2143
$wsAll.Range("F27").Value = 578
$wsAll.Range("F31").Value = 2660
$wsAll.Range("F32").Value = 478
$wsAll.Range("F38").Value = 667
$wsAll.Range("F42").Value = 921
$wsAll.Range("F43").Value = 1639
$wsAll.Range("F45").Value = 304
$wsAll.Range("F46").Value = 523
$wsAll.Range("F47").Value = 132
$wsAll.Range("F48").Value = 108
$wsAll.Range("F49").Value = 146

